$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.7128751226751149
$ws.Range("C2").Value = 0.3543759157529536
$ws.Range("D2").Value = 0.02081387122260736
$ws.Range("F2").Value = 0.4893753016434701
$ws.Range("G2").Value = 0.3327632245813135
$ws.Range("H2").Value = 0.4938077791304067
$ws.Range("I2").Value = 0.498200890656765
$ws.Range("L2").Value = 0.2924414607851986
$ws.Range("M2").Value = 0.1940569876331608
$ws.Range("N2").Value = 1.157805533820984
$ws.Range("O2").Value = 1.586271237881789
$ws.Range("B3").Value = 0.6352109149433147
$ws.Range("C3").Value = 0.3458021839002612
$ws.Range("D3").Value = 0.01856818486260181
$ws.Range("F3").Value = 0.4854483037129924
$ws.Range("G3").Value = 0.3302188605005938
$ws.Range("H3").Value = 0.4960073959640141
$ws.Range("I3").Value = 0.503346268012784
$ws.Range("L3").Value = 0.2905181278108699
$ws.Range("M3").Value = 0.1811856445322633
$ws.Range("N3").Value = 1.163573009443311
$ws.Range("O3").Value = 1.585234177931312
$ws.Range("B4").Value = 0.5874854472365598
$ws.Range("C4").Value = 0.3405163943951237
$ws.Range("D4").Value = 0.01718075915257344
$ws.Range("F4").Value = 0.4833754075076939
$ws.Range("G4").Value = 0.328925087773186
$ws.Range("H4").Value = 0.497600021986301
$ws.Range("I4").Value = 0.5068040267129597
$ws.Range("L4").Value = 0.2894929083530826
$ws.Range("M4").Value = 0.1733322675635804
$ws.Range("N4").Value = 1.167535931310873
$ws.Range("O4").Value = 1.585664395227084
$ws.Range("B5").Value = 0.5680285010659532
$ws.Range("C5").Value = 0.338357190792081
$ws.Range("D5").Value = 0.01661324912110729
$ws.Range("F5").Value = 0.482615773291819
$ws.Range("G5").Value = 0.3284653124447203
$ws.Range("H5").Value = 0.4983099428221109
$ws.Range("I5").Value = 0.5082881152425749
$ws.Range("L5").Value = 0.2891143457128678
$ws.Range("M5").Value = 0.1701446828150566
$ws.Range("N5").Value = 1.169257066467587
$ws.Range("O5").Value = 1.586107969906323
$ws.Range("B6").Value = 0.5647972237505314
$ws.Range("C6").Value = 0.3379983484327198
$ws.Range("D6").Value = 0.01651888721171701
$ws.Range("F6").Value = 0.4824947762986511
$ws.Range("G6").Value = 0.3283930389081107
$ws.Range("H6").Value = 0.498431504995402
$ws.Range("I6").Value = 0.5085390773683756
$ws.Range("L6").Value = 0.2890538565906056
$ws.Range("M6").Value = 0.1696161626616899
$ws.Range("N6").Value = 1.169549279978774
$ws.Range("O6").Value = 1.586197825999292
$ws.Range("B7").Value = 0.5872230762570041
$ws.Range("C7").Value = 0.3404872954137943
$ws.Range("D7").Value = 0.01717311406758881
$ws.Range("F7").Value = 0.483364818271049
$ws.Range("G7").Value = 0.3289186140740554
$ws.Range("H7").Value = 0.4976093495310892
$ws.Range("I7").Value = 0.5068237379218949
$ws.Range("L7").Value = 0.2894876440248169
$ws.Range("M7").Value = 0.1732892268219217
$ws.Range("N7").Value = 1.167558712850244
$ws.Range("O7").Value = 1.585669291322716
$ws.Range("B8").Value = 0.6861055181615541
$ws.Range("C8").Value = 0.3514242696606971
$ws.Range("D8").Value = 0.02004135426444265
$ws.Range("F8").Value = 0.4879510690554909
$ws.Range("G8").Value = 0.3318301657309561
$ws.Range("H8").Value = 0.4945159965070189
$ws.Range("I8").Value = 0.4999130498355484
$ws.Range("L8").Value = 0.2917460157591094
$ws.Range("M8").Value = 0.1896087759260681
$ws.Range("N8").Value = 1.15970674865909
$ws.Range("O8").Value = 1.585692178326426
$ws.Range("B9").Value = 0.8796467481263903
$ws.Range("C9").Value = 0.3726933418084286
$ws.Range("D9").Value = 0.02559688546015337
$ws.Range("F9").Value = 0.4996292743078072
$ws.Range("G9").Value = 0.3396736686073609
$ws.Range("H9").Value = 0.4903687672619981
$ws.Range("I9").Value = 0.4887309140553384
$ws.Range("L9").Value = 0.2974078541453338
$ws.Range("M9").Value = 0.2219967606496454
$ws.Range("N9").Value = 1.147647358681461
$ws.Range("O9").Value = 1.594207525287857
$ws.Range("B10").Value = 1.021557652276442
$ws.Range("C10").Value = 0.3882019363155962
$ws.Range("D10").Value = 0.02963527850752001
$ws.Range("F10").Value = 0.5098481660411665
$ws.Range("G10").Value = 0.3467435912101564
$ws.Range("H10").Value = 0.4884896626476376
$ws.Range("I10").Value = 0.4819620865078491
$ws.Range("L10").Value = 0.3023170994370332
$ws.Range("M10").Value = 0.246017793309143
$ws.Range("N10").Value = 1.140813076142351
$ws.Range("O10").Value = 1.605637508367835
$ws.Range("B11").Value = 1.086043024602475
$ws.Range("C11").Value = 0.3952297802741782
$ws.Range("D11").Value = 0.03146283184143783
$ws.Range("F11").Value = 0.5148534717729447
$ws.Range("G11").Value = 0.3502451783680414
$ws.Range("H11").Value = 0.4878880291655321
$ws.Range("I11").Value = 0.479197355579533
$ws.Range("L11").Value = 0.3047127543233614
$ws.Range("M11").Value = 0.256992722639545
$ws.Range("N11").Value = 1.138141984070728
$ws.Range("O11").Value = 1.611963070022796
$ws.Range("B12").Value = 1.110450491496692
$ws.Range("C12").Value = 0.3978869597029586
$ws.Range("D12").Value = 0.03215348040070154
$ws.Range("F12").Value = 0.5168001480095086
$ws.Range("G12").Value = 0.3516122714184178
$ws.Range("H12").Value = 0.4876965756852343
$ws.Range("I12").Value = 0.4781956788440098
$ws.Range("L12").Value = 0.3056432266194378
$ws.Range("M12").Value = 0.2611552830679287
$ws.Range("N12").Value = 1.137193321482357
$ws.Range("O12").Value = 1.614520424965889
$ws.Range("B13").Value = 1.105194454423668
$ws.Range("C13").Value = 0.3973148749182656
$ws.Range("D13").Value = 0.03200479998547223
$ws.Range("F13").Value = 0.5163786162089394
$ws.Range("G13").Value = 0.3513160136157865
$ws.Range("H13").Value = 0.4877361915200851
$ws.Range("I13").Value = 0.478409393798735
$ws.Range("L13").Value = 0.3054417979856652
$ws.Range("M13").Value = 0.260258512088221
$ws.Range("N13").Value = 1.137394841131538
$ws.Range("O13").Value = 1.613962446446351
$ws.Range("B14").Value = 1.088051286066673
$ws.Range("C14").Value = 0.3954484716845457
$ws.Range("D14").Value = 0.03151968030134356
$ws.Range("F14").Value = 0.5150125984498501
$ws.Range("G14").Value = 0.3503568254667613
$ws.Range("H14").Value = 0.4878715494238008
$ws.Range("I14").Value = 0.4791140395295024
$ws.Range("L14").Value = 0.3047888384524526
$ws.Range("M14").Value = 0.2573350483373034
$ws.Range("N14").Value = 1.138062678778979
$ws.Range("O14").Value = 1.612170217925524
$ws.Range("B15").Value = 1.077549023846871
$ws.Range("C15").Value = 0.3943047038195289
$ws.Range("D15").Value = 0.03122234629716303
$ws.Range("F15").Value = 0.5141825502581412
$ws.Range("G15").Value = 0.3497746518984002
$ws.Range("H15").Value = 0.4879591957213449
$ws.Range("I15").Value = 0.4795515519543088
$ws.Range("L15").Value = 0.3043919129583657
$ws.Range("M15").Value = 0.2555451933101125
$ws.Range("N15").Value = 1.138479925767271
$ws.Range("O15").Value = 1.611093526072295
$ws.Range("B16").Value = 1.01734182469221
$ws.Range("C16").Value = 0.3877420864715191
$ws.Range("D16").Value = 0.02951564875637303
$ws.Range("F16").Value = 0.5095282336552671
$ws.Range("G16").Value = 0.3465205042665218
$ws.Range("H16").Value = 0.4885340724648728
$ws.Range("I16").Value = 0.4821490998861542
$ws.Range("L16").Value = 0.3021638003852161
$ws.Range("M16").Value = 0.245301492486611
$ws.Range("N16").Value = 1.140996436935751
$ws.Range("O16").Value = 1.605246783400133
$ws.Range("B17").Value = 0.9803874410715707
$ws.Range("C17").Value = 0.3837090397518068
$ws.Range("D17").Value = 0.02846617848897637
$ws.Range("F17").Value = 0.5067643057885647
$ws.Range("G17").Value = 0.344597350539928
$ws.Range("H17").Value = 0.4889515615186326
$ws.Range("I17").Value = 0.4838231837393039
$ws.Range("L17").Value = 0.3008384774688437
$ws.Range("M17").Value = 0.2390293324013726
$ws.Range("N17").Value = 1.142652279088104
$ws.Range("O17").Value = 1.601948454012984
$ws.Range("B18").Value = 0.9591257058573888
$ws.Range("C18").Value = 0.3813868006922405
$ws.Range("D18").Value = 0.02786165589240142
$ws.Range("F18").Value = 0.5052081410424378
$ws.Range("G18").Value = 0.3435180665451725
$ws.Range("H18").Value = 0.4892155228439208
$ws.Range("I18").Value = 0.484815666346563
$ws.Range("L18").Value = 0.3000914779364194
$ws.Range("M18").Value = 0.2354262514345464
$ws.Range("N18").Value = 1.143645897121814
$ws.Range("O18").Value = 1.600157322348508
$ws.Range("B19").Value = 0.9519257709649196
$ws.Range("C19").Value = 0.3806001012792137
$ws.Range("D19").Value = 0.02765682228274358
$ws.Range("F19").Value = 0.5046870170891324
$ws.Range("G19").Value = 0.3431572510976224
$ws.Range("H19").Value = 0.4893089899601364
$ws.Range("I19").Value = 0.485156785133487
$ws.Range("L19").Value = 0.2998411849813039
$ws.Range("M19").Value = 0.234207090799174
$ws.Range("N19").Value = 1.143989403300999
$ws.Range("O19").Value = 1.599569075517564
$ws.Range("B20").Value = 0.9843219878404739
$ws.Range("C20").Value = 0.3841386286267721
$ws.Range("D20").Value = 0.02857798928359045
$ws.Range("F20").Value = 0.5070550560018319
$ws.Range("G20").Value = 0.3447992927774521
$ws.Range("H20").Value = 0.488904652757725
$ws.Range("I20").Value = 0.4836419114358108
$ws.Range("L20").Value = 0.3009779782506996
$ws.Range("M20").Value = 0.2396965507235294
$ws.Range("N20").Value = 1.142471746636289
$ws.Range("O20").Value = 1.602288597637227
$ws.Range("B21").Value = 1.09308698156849
$ws.Range("C21").Value = 0.395996792893925
$ws.Range("D21").Value = 0.03166221016039827
$ws.Range("F21").Value = 0.5154124394997552
$ws.Range("G21").Value = 0.3506374457854804
$ws.Range("H21").Value = 0.4878308046791062
$ws.Range("I21").Value = 0.4789058390688083
$ws.Range("L21").Value = 0.304979996939224
$ws.Range("M21").Value = 0.2581935636332346
$ws.Range("N21").Value = 1.137864814810612
$ws.Range("O21").Value = 1.612692241549524
$ws.Range("B22").Value = 1.164102001976744
$ws.Range("C22").Value = 0.4037227175220721
$ws.Range("D22").Value = 0.03366971365439753
$ws.Range("F22").Value = 0.5211733379900494
$ws.Range("G22").Value = 0.354692715343532
$ws.Range("H22").Value = 0.487340970702661
$ws.Range("I22").Value = 0.4760744149304657
$ws.Range("L22").Value = 0.3077312676188484
$ws.Range("M22").Value = 0.2703207394784641
$ws.Range("N22").Value = 1.1352200302059
$ws.Range("O22").Value = 1.620435932069881
$ws.Range("B23").Value = 1.126206805107927
$ws.Range("C23").Value = 0.3996015200701777
$ws.Range("D23").Value = 0.03259903453929525
$ws.Range("F23").Value = 0.5180712971718364
$ws.Range("G23").Value = 0.352506385995099
$ws.Range("H23").Value = 0.4875830190346733
$ws.Range("I23").Value = 0.4775614383766786
$ws.Range("L23").Value = 0.3062504659359462
$ws.Range("M23").Value = 0.2638448174673158
$ws.Range("N23").Value = 1.136598148063896
$ws.Range("O23").Value = 1.616216547822631
$ws.Range("B24").Value = 0.9825432294852021
$ws.Range("C24").Value = 0.383944422651723
$ws.Range("D24").Value = 0.02852744325883094
$ws.Range("F24").Value = 0.5069235054918053
$ws.Range("G24").Value = 0.344707912580148
$ws.Range("H24").Value = 0.4889257856305704
$ws.Range("I24").Value = 0.4837237711086431
$ws.Range("L24").Value = 0.3009148633866801
$ws.Range("M24").Value = 0.2393948923689493
$ws.Range("N24").Value = 1.142553235603152
$ws.Range("O24").Value = 1.602134491265446
$ws.Range("B25").Value = 0.8273342248923541
$ws.Range("C25").Value = 0.3669595248959041
$ws.Range("D25").Value = 0.0241014802005779
$ws.Range("F25").Value = 0.4961823615164178
$ws.Range("G25").Value = 0.3373226729840013
$ws.Range("H25").Value = 0.4912854754238936
$ws.Range("I25").Value = 0.4915020714289753
$ws.Range("L25").Value = 0.2957443917369886
$ws.Range("M25").Value = 0.2131946616254368
$ws.Range("N25").Value = 1.150553338523501
$ws.Range("O25").Value = 1.590996036994881
